$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.849.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.952.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.516"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.947.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.466"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("E13").Value = "  -3.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.834.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.439.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +14.33%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.947.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "447.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.697"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.98%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.05%  "
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("E32").Value = "  -5.86%  "
$ws.Range("E33").Value = "  +5.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.970"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.303"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.121"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("E43").Value = "  -5.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "383.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0351"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.678.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.50%  "
